# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
# Numeric-looking "Price" strings (e.g. "1.001") are written with a leading
# apostrophe so Excel stores them as literal text (matching the workbook's
# original inlineStr cells) instead of silently parsing them as numbers; the
# style is then reset to 'Normal' so no stray number-format/quote-prefix style
# sticks to the cell (the diff never touches cell styles).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "31.023.80"
$ws.Range("E2").Value = "  +1.22%  "

# Row 3
$ws.Range("D3").Value = "1.957.66"
$ws.Range("E3").Value = "  -0.18%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "'244.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.94%  "

# Row 6
$ws.Range("E6").Value = "  +0.01%  "

# Row 7
$ws.Range("D7").Value = "'0.4866"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.63%  "

# Row 8
$ws.Range("D8").Value = "'0.2954"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.34%  "

# Row 9
$ws.Range("D9").Value = "'0.06813"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.41%  "

# Row 10
$ws.Range("D10").Value = "'19.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.25%  "

# Row 11
$ws.Range("D11").Value = "'107.17"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.83%  "

# Row 12
$ws.Range("D12").Value = "1.966.12"
$ws.Range("E12").Value = "  +0.03%  "

# Row 13
$ws.Range("D13").Value = "'0.07804"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.79%  "

# Row 14
$ws.Range("D14").Value = "'5.449"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.54%  "

# Row 15
$ws.Range("D15").Value = "'0.7026"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.06%  "

# Row 16
$ws.Range("D16").Value = "'283.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.88%  "

# Row 17
$ws.Range("D17").Value = "31.033.53"
$ws.Range("E17").Value = "  +1.23%  "

# Row 18
$ws.Range("D18").Value = "'13.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.61%  "

# Row 19
$ws.Range("D19").Value = "'0.000007677"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.03%  "

# Row 20
$ws.Range("D20").Value = "2.213.97"
$ws.Range("E20").Value = "  -0.89%  "

# Row 21
$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.02%  "

# Row 22
$ws.Range("D22").Value = "'5.499"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.45%  "

# Row 24
$ws.Range("D24").Value = "'6.487"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.90%  "

# Row 25
$ws.Range("D25").Value = "'9.789"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.43%  "

# Row 26
$ws.Range("D26").Value = "'170.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.20%  "

# Row 27
$ws.Range("D27").Value = "'19.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.67%  "

# Row 28
$ws.Range("D28").Value = "'2.213"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.59%  "

# Row 29
$ws.Range("D29").Value = "'0.1058"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.65%  "

# Row 30
$ws.Range("D30").Value = "'1.410"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.95%  "

# Row 31
$ws.Range("E31").Value = "  -1.95%  "

# Row 32
$ws.Range("D32").Value = "'4.598"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.48%  "

# Row 33
$ws.Range("D33").Value = "'4.453"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.12%  "

# Row 34
$ws.Range("D34").Value = "'0.04927"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.75%  "

# Row 35
$ws.Range("D35").Value = "'0.7599"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.53%  "

# Row 36
$ws.Range("E36").Value = "  -0.83%  "

# Row 37
$ws.Range("D37").Value = "'2.729"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.02%  "

# Row 38
$ws.Range("D38").Value = "'0.02012"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.20%  "

# Row 39
$ws.Range("D39").Value = "'2.702"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.36%  "

# Row 40
$ws.Range("D40").Value = "'6.519"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.12%  "

# Row 41
$ws.Range("D41").Value = "'2.112"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.98%  "

# Row 42
$ws.Range("D42").Value = "'76.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.93%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.8874"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.58%  "

# Row 44
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.4462"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.50%  "

# Row 45
$ws.Range("D45").Value = "'109.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.25%  "

# Row 46
$ws.Range("D46").Value = "'8.158"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.25%  "

# Row 47
$ws.Range("D47").Value = "'1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.01%  "

# Row 48
$ws.Range("D48").Value = "'996.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.50%  "

# Row 49
$ws.Range("D49").Value = "'0.1259"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.03%  "

# Row 50
$ws.Range("D50").Value = "'9.319"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.92%  "

# Row 51
$ws.Range("E51").Value = "  -0.64%  "
